$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, pushing the existing rows 19..129 down to 20..130
$ws.Rows("19:19").Insert()

# Populate the freshly inserted row 19 with the new weekly price-report record
# (same market/category metadata as the record that used to sit there, but a
# new date and new volume/price figures).
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = "Femacal de La Calera"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 44677
$ws.Range("E19").Value = 5
$ws.Range("F19").Value = 100112026
$ws.Range("G19").Value = "Haba"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 45
$ws.Range("K19").Value = 21000
$ws.Range("L19").Value = 21000
$ws.Range("M19").Value = 21000
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Provincia de Limarí"
$ws.Range("P19").Value = 840
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
